# Upload most of the available VKG resources from ontop-examples.
#
# Row 2 ("uobm") is replaced by "bsbm", and five new rows are appended
# (npd, cordis, suedtirol, canonical, dblp) mirroring the VKG resource
# table that ships with ontop-examples.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 already carries the "hyperlink" cell style (s="1"). Stash a copy of
# that formatting on a scratch cell far outside the used range so we can
# re-apply it later to every hyperlinked cell in column B, regardless of
# what Hyperlinks.Add() below does to the cell's style.
$ws.Range("B2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 2: bsbm (replaces uobm) ---
$ws.Range("A2").Value = "bsbm"
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B2").Value = "https://github.com/ontop/ontop-examples/tree/master/dke-2022-mapping-patterns/scenarios"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/ontop/ontop-examples/tree/master/dke-2022-mapping-patterns/scenarios")
$ws.Range("C2").Value = "√"
$ws.Range("D2").Value = "√"
$ws.Range("E2").Value = "√"
$ws.Range("F2").Value = "×"
$ws.Range("G2").Value = "MySQL"

# --- Row 3: npd ---
$ws.Range("A3").Value = "npd"
$ws.Range("B3").Value = "https://github.com/ontop/ontop-examples/tree/master/caise-2021-patterns/scenarios/npd"
$ws.Range("C3").Value = "√"
$ws.Range("D3").Value = "√"
$ws.Range("E3").Value = "√"
$ws.Range("F3").Value = "√"
$ws.Range("G3").Value = "MySQL"

# --- Row 4: cordis ---
$ws.Range("A4").Value = "cordis"
$ws.Range("B4").Value = "https://github.com/ontop/ontop-examples/tree/master/caise-2021-patterns/scenarios/cordis"
$ws.Range("C4").Value = "√"
$ws.Range("D4").Value = "√"
$ws.Range("E4").Value = "√"
$ws.Range("F4").Value = "√"
$ws.Range("G4").Value = "PostgreSQL"

# --- Row 5: suedtirol ---
$ws.Range("A5").Value = "suedtirol"
$ws.Range("B5").Value = "https://github.com/ontop/ontop-examples/tree/master/dke-2022-mapping-patterns/scenarios/suedtirol-open-data"
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/ontop/ontop-examples/tree/master/dke-2022-mapping-patterns/scenarios/suedtirol-open-data")
$ws.Range("C5").Value = "√"
$ws.Range("D5").Value = "√"
$ws.Range("E5").Value = "√"
$ws.Range("F5").Value = "√"
$ws.Range("G5").Value = "PostgreSQL"

# --- Row 6: canonical ---
$ws.Range("A6").Value = "canonical"
$ws.Range("B6").Value = "https://github.com/ontop/ontop-examples/tree/master/eswc-2018-canonical-iri"
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/ontop/ontop-examples/tree/master/eswc-2018-canonical-iri")
$ws.Range("C6").Value = "√"
$ws.Range("D6").Value = "√"
$ws.Range("E6").Value = "√"
$ws.Range("F6").Value = "×"
$ws.Range("G6").Value = "PostgreSQL"

# --- Row 7: dblp ---
$ws.Range("A7").Value = "dblp"
$ws.Range("B7").Value = "https://github.com/ontop/ontop-examples/tree/master/swj-2017-vig"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/ontop/ontop-examples/tree/master/swj-2017-vig")
$ws.Range("C7").Value = "√"
$ws.Range("D7").Value = "√"
$ws.Range("E7").Value = "√"
$ws.Range("F7").Value = "×"
$ws.Range("G7").Value = "MySQL"

# Hyperlinks.Add() resets the cell style on every hyperlinked B cell -
# restore the original hyperlink look (style s="1") on all of them.
$ws.Range("Z100").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# Match the final selection recorded in the sheet view.
$ws.Range("F16").Select()
